$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.259.91"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "1.860.06"
$ws.Range("E3").Value = "  -0.54%  "
$origStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("D4").Style = $origStyle
$ws.Range("E4").Value = "  +0.21%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.80"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  -0.29%  "
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6955"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  -1.84%  "
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9996"
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = "  +0.13%  "
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07801"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  -0.57%  "
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3078"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  -2.43%  "
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.90"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  -2.19%  "
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07824"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  -2.33%  "
$ws.Range("D12").Value = "1.858.89"
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "92.55"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  -2.09%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.117"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  -2.23%  "
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6904"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  -2.05%  "
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.549"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  +1.83%  "
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008429"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  +1.54%  "
$ws.Range("D18").Value = "29.252.44"
$ws.Range("E18").Value = "  -0.68%  "
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "248.03"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  -3.51%  "
$ws.Range("D20").Value = "2.110.59"
$ws.Range("E20").Value = "  -0.80%  "
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.85"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  -2.86%  "
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  +0.15%  "
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.568"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  -1.08%  "
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9996"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  +0.18%  "
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1501"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  -4.47%  "
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.87"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  +0.95%  "
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.884"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  -1.81%  "
$ws.Range("E28").Value = "  -2.01%  "
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.556"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  +4.01%  "
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.270"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  -1.59%  "
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.219"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  -1.40%  "
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.200"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  -1.37%  "
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05220"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  -1.72%  "
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7624"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  +1.75%  "
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.851"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  -2.62%  "
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.172"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  +0.46%  "
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.707"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E38").Value = "  -0.90%  "
$ws.Range("D39").Value = "1.234.69"
$ws.Range("E39").Value = "  -1.32%  "
$ws.Range("E40").Value = "  -0.71%  "
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9102"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  +1.22%  "
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.69"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  -0.59%  "
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9990"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  +0.22%  "
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.556"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  -9.54%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000124"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  -5.11%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "2.008.67"
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.601"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  +0.77%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "65.22"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  -8.46%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5182"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  +0.07%  "
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.753"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  -2.11%  "
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.041"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  -0.32%  "
